$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the "2021-Q4" sheet to create "2022-Q1" (same column layout /
#    header / styles), placed right after "2021-Q4" and before "总计".
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2021-Q4")
$ws1.Copy($null, $ws1)
$wsNew = $wb.Worksheets.Item("2021-Q4 (2)")
$wsNew.Name = "2022-Q1"
$ws2 = $wb.Worksheets.Item("2022-Q1")

# Overwrite the fund rows with the 2022-Q1 figures (headers / A-column /
# styles were already copied verbatim from "2021-Q4").
$ws2.Range("B2").Value = "'012348"
$ws2.Range("C2").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）A"
$ws2.Range("D2").Value = "'38.10"
$ws2.Range("E2").Value = "'92.34"
$ws2.Range("F2").Value = "'7.95"
$ws2.Range("G2").Value = "'3.0290"
$ws2.Range("H2").Value = 2

$ws2.Range("B3").Value = "'012349"
$ws2.Range("C3").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）C"
$ws2.Range("D3").Value = "'14.77"
$ws2.Range("E3").Value = "'92.34"
$ws2.Range("F3").Value = "'7.95"
$ws2.Range("G3").Value = "'1.1742"
$ws2.Range("H3").Value = 2

# ---------------------------------------------------------------------------
# 2) "总计" sheet: insert the new 2022-Q1 summary row above the existing
#    2021-Q4 row (push it from row 2 down to row 3).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("总计")

# Clone A2's format down into A3 first, so the shifted 2021-Q4 row keeps the
# same look (bold/border/centered) as before.
$ws3.Range("A2").Copy()
$ws3.Range("A3").PasteSpecial(-4122)

# Move the old 2021-Q4 totals down to row 3.
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "2021-Q4"
$ws3.Range("C3").Value = 2
$ws3.Range("D3").Value = 1.2

# Write the new 2022-Q1 totals into row 2.
$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "2022-Q1"
$ws3.Range("C2").Value = 2
$ws3.Range("D2").Value = 4.2

# Restore the originally-active sheet/tab.
$ws1.Activate()
